$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C14").Value = "https://media.discordapp.net/attachments/1162451241872412902/1169227786918043698/land.png?ex=6554a33b&is=65422e3b&hm=1c6018c123d4f2895d36861b46ffc0c6f633e0f943f7f9c33cb747b37625ac04&=&width=380&height=380"
$ws.Range("C13").Value = "https://media.discordapp.net/attachments/1162451241872412902/1169227787199053894/kamlesh_kaniyal.png?ex=6554a33b&is=65422e3b&hm=2340ce8c17d03dee81d7ed6beb788e008fca7d12ffb7a3b53935da73dac4b1ff&=&width=380&height=380"
$ws.Range("A13").Value = "Kamlesh Kaniyal"
$ws.Range("A14").Value = "Rahul Rauleta"

$ws.Range("B13").Value = "Cricket Co-ordinator"
$ws.Range("B14").Value = "Volleyball Co-ordinator"

$ws.Range("A12").Copy()
$ws.Range("A13:A14").PasteSpecial(-4122)

$ws.Range("B12").Copy()
$ws.Range("B13:B14").PasteSpecial(-4122)

$ws.Range("G13").Select()
